$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

# Insert a new row above the current row 3 (ESBCONF), shifting everything
# below down by one. This mirrors Excel's "Insert Sheet Rows" behaviour.
$ws.Rows.Item(3).Insert()

# The row above (TESTDB, row 2) carries the "Good" highlight style used for
# every populated row in this block; apply the same style to the cells the
# new row will actually use, matching the sheet's existing pattern.
$ws.Range("A3:D3").Style = "Good"
$ws.Range("F3").Style = "Good"
$ws.Range("H3").Style = "Good"

# Populate the newly inserted row 3 with the RESOURCELINK data source.
$ws.Range("A3").Value = "RESOURCELINK"
$ws.Range("B3").Value = "CMTEST"
$ws.Range("C3").Value = "jm08_cmt"
$ws.Range("D3").Formula = '=CONCATENATE( "mqsisetdbparms ",ConfigData!$D$4," -n ",A3," -u ",B3," -p ",C3)'
$ws.Range("F3").Formula = '=CONCATENATE( "mqsicvp ",ConfigData!$D$4," -n ",A3)'
$ws.Range("H3").Value = "Attempt to standardise the DSN across environments."

# Match the author's final selection.
$ws.Range("F3").Select()
